$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.157.58"
$ws.Range("E2").Value = "  -2.24%  "

$ws.Range("D3").Value = "1.839.06"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.04"
$ws.Range("E5").Value = "  -2.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6824"
$ws.Range("E6").Value = "  -2.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3000"
$ws.Range("E8").Value = "  -2.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07468"
$ws.Range("E9").Value = "  -4.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.25"
$ws.Range("E10").Value = "  -2.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07643"

$ws.Range("D12").Value = "1.837.74"
$ws.Range("E12").Value = "  -1.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.042"
$ws.Range("E13").Value = "  -2.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6814"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.99"
$ws.Range("E15").Value = "  -5.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.126"
$ws.Range("E16").Value = "  -8.24%  "

$ws.Range("D17").Value = "29.152.22"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008229"
$ws.Range("E18").Value = "  -2.32%  "

$ws.Range("D19").Value = "2.084.35"
$ws.Range("E19").Value = "  -1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.27"
$ws.Range("E20").Value = "  -5.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("E21").Value = "  -2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.349"
$ws.Range("E23").Value = "  -4.16%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.76"
$ws.Range("E25").Value = "  +0.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1437"
$ws.Range("E26").Value = "  -5.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.713"
$ws.Range("E27").Value = "  -3.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  -1.95%  "

$ws.Range("E29").Value = "  -2.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.266"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05386"
$ws.Range("E33").Value = "  +5.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7561"
$ws.Range("E34").Value = "  -4.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.858"
$ws.Range("E35").Value = "  -4.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.134"
$ws.Range("E36").Value = "  -3.07%  "

$ws.Range("E37").Value = "  -0.90%  "

$ws.Range("D38").Value = "1.312.70"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01831"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.724"
$ws.Range("E40").Value = "  -1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9472"
$ws.Range("E41").Value = "  -2.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.022"
$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.63"
$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").Value = "1.985.16"
$ws.Range("E45").Value = "  -1.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5177"
$ws.Range("E46").Value = "  -0.55%  "

$ws.Range("E47").Value = "  -3.52%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.49"
$ws.Range("E48").Value = "  -1.76%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.774"
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("B50").Value = "XinFinNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07685"
$ws.Range("E50").Value = "  +15.52%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.421"
$ws.Range("E51").Value = "  -4.24%  "
